$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 53287
$ws.Cells.Item(2, 2).Value = "Maria Clara Rezende"
$ws.Cells.Item(2, 3).Value = "Recursos Humanos"
$ws.Cells.Item(2, 4).Value = "Outros"
$ws.Cells.Item(2, 5).Value = 6
$ws.Cells.Item(2, 6).Value = 45083
$ws.Cells.Item(2, 7).Value = 6966.97

# Row 3
$ws.Cells.Item(3, 1).Value = 96045
$ws.Cells.Item(3, 2).Value = "Dr. Noah Costa"
$ws.Cells.Item(3, 3).Value = "Juridico"
$ws.Cells.Item(3, 4).Value = "Viagem de negocios"
$ws.Cells.Item(3, 5).Value = 4
$ws.Cells.Item(3, 6).Value = 45078
$ws.Cells.Item(3, 7).Value = 4343.99

# Row 4
$ws.Cells.Item(4, 1).Value = 79748
$ws.Cells.Item(4, 2).Value = "Danilo Guerra"
$ws.Cells.Item(4, 3).Value = "Juridico"
$ws.Cells.Item(4, 4).Value = "Viagem de negocios"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 45086
$ws.Cells.Item(4, 7).Value = 9030.110000000001

# Row 5
$ws.Cells.Item(5, 1).Value = 38307
$ws.Cells.Item(5, 2).Value = "Antônio Montenegro"
$ws.Cells.Item(5, 3).Value = "TI"
$ws.Cells.Item(5, 4).Value = "Outros"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 45104
$ws.Cells.Item(5, 7).Value = 2816.19

# Row 6
$ws.Cells.Item(6, 1).Value = 83757
$ws.Cells.Item(6, 2).Value = "Pietro da Rosa"
$ws.Cells.Item(6, 3).Value = "Vendas"
$ws.Cells.Item(6, 4).Value = "Viagem de negocios"
$ws.Cells.Item(6, 5).Value = 7
$ws.Cells.Item(6, 6).Value = 45105
$ws.Cells.Item(6, 7).Value = 6620.87

# Row 7
$ws.Cells.Item(7, 1).Value = 72727
$ws.Cells.Item(7, 2).Value = "Srta. Ester Teixeira"
$ws.Cells.Item(7, 3).Value = "Recursos Humanos"
$ws.Cells.Item(7, 4).Value = "Viagem de negocios"
$ws.Cells.Item(7, 5).Value = 7
$ws.Cells.Item(7, 6).Value = 45096
$ws.Cells.Item(7, 7).Value = 8369.620000000001

# Row 8
$ws.Cells.Item(8, 1).Value = 4302
$ws.Cells.Item(8, 2).Value = "Heitor Rocha"
$ws.Cells.Item(8, 3).Value = "Juridico"
$ws.Cells.Item(8, 4).Value = "Viagem de negocios"
$ws.Cells.Item(8, 5).Value = 4
$ws.Cells.Item(8, 6).Value = 45105
$ws.Cells.Item(8, 7).Value = 4518.93

# Row 9
$ws.Cells.Item(9, 1).Value = 39378
$ws.Cells.Item(9, 2).Value = "Rodrigo Casa Grande"
$ws.Cells.Item(9, 3).Value = "Recursos Humanos"
$ws.Cells.Item(9, 4).Value = "Outros"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 45083
$ws.Cells.Item(9, 7).Value = 6837.33

# Row 10
$ws.Cells.Item(10, 1).Value = 14991
$ws.Cells.Item(10, 2).Value = "Brenda Pires"
$ws.Cells.Item(10, 3).Value = "Financeiro"
$ws.Cells.Item(10, 4).Value = "Problemas pessoais"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 45103
$ws.Cells.Item(10, 7).Value = 5324.24

# Row 11
$ws.Cells.Item(11, 1).Value = 16432
$ws.Cells.Item(11, 2).Value = "João Felipe Sousa"
$ws.Cells.Item(11, 3).Value = "Marketing"
$ws.Cells.Item(11, 4).Value = "Consulta medica"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 45079
$ws.Cells.Item(11, 7).Value = 3591.24
